$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date in C1 (45294 -> 45379) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: update priority values for the "hard coal" row (row 3) ---
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# --- Update the active selection on the FPIEBP sheet (F4 -> E3) ---
$wsFPIEBP.Activate()
$wsFPIEBP.Range("E3").Select()
